$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "328.11"
Set-TextValue $ws.Range("E2") "-0.25%"

Set-TextValue $ws.Range("D3") "44.29"
Set-TextValue $ws.Range("E3") "-0.09%"

Set-TextValue $ws.Range("D4") "5.107"
Set-TextValue $ws.Range("E4") "-7.40%"

Set-TextValue $ws.Range("D5") "0.08388"
Set-TextValue $ws.Range("E5") "3.92%"

Set-TextValue $ws.Range("B6") "GateToken"
Set-TextValue $ws.Range("C6") "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D6") "4.447"
Set-TextValue $ws.Range("E6") "0.27%"

Set-TextValue $ws.Range("B7") "FTXToken"
Set-TextValue $ws.Range("C7") "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D7") "1.945"
Set-TextValue $ws.Range("E7") "-5.40%"

Set-TextValue $ws.Range("B8") "MXToken"
Set-TextValue $ws.Range("C8") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D8") "0.9748"
Set-TextValue $ws.Range("E8") "1.57%"

Set-TextValue $ws.Range("B9") "BTSEToken"
Set-TextValue $ws.Range("C9") "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D9") "2.500"
Set-TextValue $ws.Range("E9") "-4.62%"

Set-TextValue $ws.Range("B10") "LiechtensteinCryptoassetsExchange"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D10") "0.1140"
Set-TextValue $ws.Range("E10") "0.77%"

Set-TextValue $ws.Range("B11") "WazirX"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D11") "0.1903"
Set-TextValue $ws.Range("E11") "1.19%"

Set-TextValue $ws.Range("B12") "MandalaExchangeToken"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D12") "0.09704"
Set-TextValue $ws.Range("E12") "-2.58%"

Set-TextValue $ws.Range("B13") "BitrueCoin"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D13") "0.04611"
Set-TextValue $ws.Range("E13") "-1.91%"

Set-TextValue $ws.Range("B14") "BitMartToken"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D14") "0.1060"
Set-TextValue $ws.Range("E14") "-0.05%"

Set-TextValue $ws.Range("B15") "BitForexToken"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D15") "0.001281"
Set-TextValue $ws.Range("E15") "1.71%"

Set-TextValue $ws.Range("B16") "TigerCash"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D16") "0.005895"
Set-TextValue $ws.Range("E16") "-3.83%"

Set-TextValue $ws.Range("B17") "LEO"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D17") "3.403"
Set-TextValue $ws.Range("E17") "1.89%"

Set-TextValue $ws.Range("D18") "0.3362"
Set-TextValue $ws.Range("E18") "1.41%"

Set-TextValue $ws.Range("D19") "8.848"
Set-TextValue $ws.Range("E19") "-13.44%"

Set-TextValue $ws.Range("D20") "0.1364"
Set-TextValue $ws.Range("E20") "-2.22%"

Set-TextValue $ws.Range("D21") "0.2581"
Set-TextValue $ws.Range("E21") "0.04%"

Set-TextValue $ws.Range("D22") "0.04169"
Set-TextValue $ws.Range("E22") "1.64%"

Set-TextValue $ws.Range("D23") "0.001253"
Set-TextValue $ws.Range("E23") "-4.55%"

Set-TextValue $ws.Range("D24") "0.004426"
Set-TextValue $ws.Range("E24") "1.61%"

Set-TextValue $ws.Range("D25") "0.0001303"
Set-TextValue $ws.Range("E25") "1.47%"

Set-TextValue $ws.Range("E26") "-20.37%"

Set-TextValue $ws.Range("D38") "0.02739"
Set-TextValue $ws.Range("E38") "3.77%"

Set-TextValue $ws.Range("D39") "0.05627"
Set-TextValue $ws.Range("E39") "0.09%"

Set-TextValue $ws.Range("D40") "0.007882"
Set-TextValue $ws.Range("E40") "3.59%"

Set-TextValue $ws.Range("E41") "0.69%"

Set-TextValue $ws.Range("E42") "-0.48%"

Set-TextValue $ws.Range("D43") "0.002114"
Set-TextValue $ws.Range("E43") "6.13%"

Set-TextValue $ws.Range("D44") "0.007916"
Set-TextValue $ws.Range("E44") "-9.40%"

Set-TextValue $ws.Range("D45") "0.3505"

Set-TextValue $ws.Range("D46") "0.00006924"
Set-TextValue $ws.Range("E46") "-2.74%"

Set-TextValue $ws.Range("E47") "-0.03%"

Set-TextValue $ws.Range("D48") "0.003505"
Set-TextValue $ws.Range("E48") "0.06%"

Set-TextValue $ws.Range("D49") "0.003540"
Set-TextValue $ws.Range("E49") "40.03%"

Set-TextValue $ws.Range("E50") "-0.03%"

Set-TextValue $ws.Range("D51") "0.0002006"
Set-TextValue $ws.Range("E51") "-0.03%"
